$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Delete Event": the "TestEvent" row (row 4) is removed via the
# spreadsheet's delete-event routine, so every row below it shifts up one.
$ws.Rows(4).Delete()

# Format the two affected rows as Text first so the big numeric-looking
# IDs/day/month/year values are stored as text (no scientific notation,
# no lost leading zeros) instead of being auto-coerced to numbers.
$rng = $ws.Range("A5:F6")
$rng.NumberFormat = "@"

# Row 5 (the shifted-up, previously blank template row) becomes a newly
# added event.
$ws.Cells.Item(5, 1).Value = "8090963507605344749"
$ws.Cells.Item(5, 2).Value = "11"
$ws.Cells.Item(5, 3).Value = "12"
$ws.Cells.Item(5, 4).Value = "2022"
$ws.Cells.Item(5, 6).Value = "5"

# Row 6: another newly added event.
$ws.Cells.Item(6, 1).Value = "9527163690503573997"
$ws.Cells.Item(6, 2).Value = "17"
$ws.Cells.Item(6, 3).Value = "10"
$ws.Cells.Item(6, 4).Value = "2022"
$ws.Cells.Item(6, 5).Value = "Test"
$ws.Cells.Item(6, 6).Value = "0"

# Drop back to the default "Normal" style so the cells don't retain an
# explicit Text-format style index (keeps the stored type as text though).
$rng.Style = "Normal"

# Row 5 never got a description typed in, so that cell stays blank/empty.
$ws.Cells.Item(5, 5).ClearContents()

$ws.Range("N13").Select()
